$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheet "Third" -> "NamedRange" and populate it with the named
#    range demo content (labels, sum formulas, boxed data, merged cell).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Third")
$ws.Name = "NamedRange"

# Clear previous placeholder content (B1/D1/F1/H1/A3/C3/E3/G3/B5/D5/F5/H5).
$ws.Cells.Clear()

# Column widths used by the new layout.
$ws.Columns.Item(1).ColumnWidth = 15.5
$ws.Columns.Item(2).ColumnWidth = 9.75
$ws.Columns.Item(6).ColumnWidth = 11.75

# Row 1: headers + sum formulas.
$ws.Rows.Item(1).RowHeight = 21

$ws.Range("A1").Value = "named range"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 16

$ws.Range("B1").Value = "TestRange1"
$ws.Range("B1").Font.Bold = $true

$ws.Range("F1").Value = "RangeMerged"
$ws.Range("F1").Font.Bold = $true

# Data boxed by TestRange1 (B2:D3).
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6
$ws.Range("B2:D3").BorderAround(1)

# Merged cell RangeMerged (F2:H3) plus the F4:H4 helper row.
$ws.Range("F2:H3").Merge()
$ws.Range("F2").Value = 1
$ws.Range("F2:H3").BorderAround(1)

$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 4

# ---------------------------------------------------------------------------
# 2) Re-point the defined names at the new local sheet/range instead of the
#    external workbook. Deleting + re-adding (instead of just assigning
#    .RefersTo) avoids a stale "external reference" flag on the name.
# ---------------------------------------------------------------------------
$wb.Names.Item("RangeMerged").Delete()
$wb.Names.Add("RangeMerged", "=NamedRange!`$F`$2")
$wb.Names.Item("TestRange1").Delete()
$wb.Names.Add("TestRange1", "=NamedRange!`$B`$2:`$D`$3")

# Now that the names resolve locally, add the SUM formulas that reference them.
$ws.Range("C1").Formula = "=SUM(TestRange1)"
$ws.Range("G1").Formula = "=SUM(RangeMerged,F4:H4)"

# ---------------------------------------------------------------------------
# 3) Update the "Format" sheet selection (unrelated cursor move captured by
#    the diff) without changing which sheet/tab is active.
# ---------------------------------------------------------------------------
$wsFormat = $wb.Worksheets.Item("Format")
$wsFormat.Range("A7").Select()

# ---------------------------------------------------------------------------
# 4) Make "NamedRange" the active sheet/tab (was "cell-border" before) and
#    leave the cursor on C7, matching the new bookViews/activeTab + the
#    sheet's own tabSelected/selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C7").Select()
